# Updated requirement specifications (Web UI sheet): the "Jobs" web-UI
# requirement section is fleshed out with a filter/data-grid layout,
# replacing the previous placeholder text and the leftover
# Ward/Category/Device grid that used to live under it. Everything that
# used to sit below that block (Medical Attendant, Associate Operators,
# Reports, Patients ...) shifts down by four rows to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web UI")

# Make room for the new "Jobs" detail rows: insert 4 blank rows right
# above the old "Medical Attendant" row (row 54), pushing it (and
# everything after it) down to row 58.
$ws.Rows("54:57").Insert()

# The rows that stayed in place (46-53) still hold the old placeholder
# text ("Need to create this page in the product", the Ward/Category/
# Device grid, etc.) - wipe them before writing the new content.
$ws.Range("B46:I54").ClearContents()

# --- New "Jobs" requirement block ---
$ws.Range("B46").Value = "Jobs"
$ws.Range("C46").Value = "Filter"
$ws.Range("C47").Value = "Vehicle No"
$ws.Range("C48").Value = "Status"

$ws.Range("C51").Value = "Service Date"
$ws.Range("D51").Value = "Today"
$ws.Range("E51").Value = "All"
$ws.Range("F51").Value = "Select"
$ws.Range("G51").Value = "Similar to graph above the datagrid"

$ws.Range("C53").Value = "Data grid columns"

$ws.Range("C54").Value = "Vehicle No"
$ws.Range("D54").Value = "Service Date"
$ws.Range("E54").Value = "In time"
$ws.Range("F54").Value = "Out Time"
$ws.Range("G54").Value = "Token"
$ws.Range("H54").Value = "Status"
$ws.Range("I54").Value = "Action"

# Leftover orphan cell from the old content, preserved at its shifted spot.
$ws.Range("E57").Value = "ICU"

# Reset the view to match where the author was working.
$ws.Activate()
$ws.Range("D54").Select()
